# Desarrollo de Cotizacion de Equipos Parte 0.4
#
# Three closing paragraphs (the "${nvchAutor}" signature line, the
# "${nvchCargo}" line and the "Resteco S.A" line) move from justified
# to centered text with a large right indent, and the stray leading
# whitespace runs in front of "${nvchAutor}" and "Resteco S.A" are
# removed. The "_GoBack" bookmark that used to wrap "Resteco S.A" is
# relocated so it starts right before "${nvchAutor}" instead (and still
# ends right after "S.A").

$d = $word.ActiveDocument

# Locate the three target paragraphs by their distinctive text so the
# script does not depend on a hard-coded paragraph index.
$pAutorIdx = -1
$pCargoIdx = -1
$pRestecoIdx = -1
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $t = $d.Paragraphs($i).Range.Text
    if ($t -match "nvchAutor") { $pAutorIdx = $i }
    elseif ($t -match "nvchCargo") { $pCargoIdx = $i }
    elseif ($t -match "Resteco") { $pRestecoIdx = $i }
}

# --- Remove the leading run of 10 spaces before "${nvchAutor}" -------
# (earlier in the document, so do it first while offsets are "clean")
$pA = $d.Paragraphs($pAutorIdx)
$aText = $pA.Range.Text
$aLeadLen = $aText.IndexOf("`${")
if ($aLeadLen -gt 0) {
    $aStart = $pA.Range.Start
    $d.Range($aStart, $aStart + $aLeadLen).Delete()
}

# --- Remove the leading run of 13 spaces before "Resteco" ------------
# Re-fetch the paragraph fresh so the Start/End reflect the edit above.
$pR = $d.Paragraphs($pRestecoIdx)
$rText = $pR.Range.Text
$rLeadLen = $rText.IndexOf("Resteco")
if ($rLeadLen -gt 0) {
    $rStart = $pR.Range.Start
    $d.Range($rStart, $rStart + $rLeadLen).Delete()
}

# --- Move the "_GoBack" bookmark from "Resteco S.A" to just before ---
# --- "${nvchAutor}", keeping its original end point ("...S.A")       -
$bm = $d.Bookmarks.Item("_GoBack")
$bmEnd = $bm.End
$bm.Delete()

$pA = $d.Paragraphs($pAutorIdx)
$newBmStart = $pA.Range.Start
$d.Bookmarks.Add("_GoBack", $d.Range($newBmStart, $bmEnd))

# --- Re-justify the three paragraphs: centered, with a big right -----
# --- indent (5812 twips = 290.6 points) instead of fully justified.  -
foreach ($idx in @($pAutorIdx, $pCargoIdx, $pRestecoIdx)) {
    $p = $d.Paragraphs($idx)
    $p.Format.RightIndent = 290.6
    $p.Format.Alignment = 1
}
